$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4869.88
$ws.Range("J64").Value = 8424.286
$ws.Range("L64").Value = 8424.286
$ws.Range("N64").Value = -8920.286
$ws.Range("H67").Value = 4869.88
$ws.Range("J67").Value = 8424.286
$ws.Range("L67").Value = 8424.286
$ws.Range("N67").Value = -10140.286
$ws.Range("H74").Value = 3428.524
$ws.Range("J74").Value = 3500
$ws.Range("L74").Value = 3500
$ws.Range("N74").Value = -5372
$ws.Range("H77").Value = 3428.524
$ws.Range("J77").Value = 3500
$ws.Range("L77").Value = 17500
$ws.Range("N77").Value = -26860
$ws.Range("H137").Value = 16668884
$ws.Range("I137").Value = 1500
$ws.Range("J137").Value = 28574158
$ws.Range("K137").Value = 4500
$ws.Range("L137").Value = 85722474
$ws.Range("M137").Value = -1950
$ws.Range("N137").Value = -85727574

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 20003350
$ws.Range("I61").Value = 29415014
$ws.Range("K61").Value = 29415014
$ws.Range("M61").Value = -29414802
$ws.Range("H74").Value = 15628274
$ws.Range("I74").Value = 23811386
$ws.Range("J74").Value = 5970.364
$ws.Range("K74").Value = 23811386
$ws.Range("L74").Value = 5970.364
$ws.Range("M74").Value = -23810512
$ws.Range("N74").Value = -7718.364
$ws.Range("H77").Value = 15628274
$ws.Range("I77").Value = 23811386
$ws.Range("J77").Value = 5970.364
$ws.Range("K77").Value = 119056930
$ws.Range("L77").Value = 29851.82
$ws.Range("M77").Value = -119052562
$ws.Range("N77").Value = -38587.82
$ws.Range("H94").Value = 50330
$ws.Range("J94").Value = 50330
$ws.Range("L94").Value = 50330
$ws.Range("N94").Value = -52132
$ws.Range("H132").Value = 35719124
$ws.Range("I132").Value = 125005010
$ws.Range("K132").Value = 375015030
$ws.Range("M132").Value = -375012500
$ws.Range("H136").Value = 20003350
$ws.Range("I136").Value = 29415014
$ws.Range("K136").Value = 88245042
$ws.Range("M136").Value = -88242492

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4546.303
$ws.Range("I105").Value = 3352
$ws.Range("J105").Value = 4928.48
$ws.Range("K105").Value = 3352
$ws.Range("L105").Value = 4928.48
$ws.Range("M105").Value = -1605
$ws.Range("N105").Value = -8422.48
$ws.Range("H134").Value = 3414.724
$ws.Range("I134").Value = 1852.0588
$ws.Range("J134").Value = 5628.5
$ws.Range("K134").Value = 5556.1764
$ws.Range("L134").Value = 16885.5
$ws.Range("M134").Value = -3021.1764
$ws.Range("N134").Value = -21955.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3754.2307
$ws.Range("I58").Value = 1750
$ws.Range("J58").Value = 6961
$ws.Range("K58").Value = 1750
$ws.Range("L58").Value = 6961
$ws.Range("M58").Value = -1547
$ws.Range("N58").Value = -7367
$ws.Range("H62").Value = 2270.2856
$ws.Range("I62").Value = 2342.5
$ws.Range("J62").Value = 2039.2
$ws.Range("K62").Value = 2342.5
$ws.Range("L62").Value = 2039.2
$ws.Range("M62").Value = -1718.5
$ws.Range("N62").Value = -3287.2
$ws.Range("H65").Value = 2270.2856
$ws.Range("I65").Value = 2342.5
$ws.Range("J65").Value = 2039.2
$ws.Range("K65").Value = 11712.5
$ws.Range("L65").Value = 10196
$ws.Range("M65").Value = -8592.5
$ws.Range("N65").Value = -16436
$ws.Range("H132").Value = 11232.286
$ws.Range("I132").Value = 14104.444
$ws.Range("J132").Value = 6062.4
$ws.Range("K132").Value = 42313.33199999999
$ws.Range("L132").Value = 18187.2
$ws.Range("M132").Value = -39783.33199999999
$ws.Range("N132").Value = -23247.2
$ws.Range("H134").Value = 2944.8572
$ws.Range("I134").Value = 3150
$ws.Range("J134").Value = 2671.3333
$ws.Range("K134").Value = 9450
$ws.Range("L134").Value = 8013.999899999999
$ws.Range("M134").Value = -6915
$ws.Range("N134").Value = -13083.9999
$ws.Range("H136").Value = 3754.2307
$ws.Range("I136").Value = 1750
$ws.Range("J136").Value = 6961
$ws.Range("K136").Value = 5250
$ws.Range("L136").Value = 20883
$ws.Range("M136").Value = -2700
$ws.Range("N136").Value = -25983

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 545.3333
$ws.Range("I5").Value = 479.23077
$ws.Range("J5").Value = 975
$ws.Range("K5").Value = 1437.69231
$ws.Range("L5").Value = 2925
$ws.Range("M5").Value = -1325.69231
$ws.Range("N5").Value = -3149
$ws.Range("H118").Value = 2563
$ws.Range("J118").Value = 2743.3333
$ws.Range("L118").Value = 8229.999899999999
$ws.Range("N118").Value = -10715.9999
$ws.Range("H135").Value = 545.3333
$ws.Range("I135").Value = 479.23077
$ws.Range("J135").Value = 975
$ws.Range("K135").Value = 4313.07693
$ws.Range("L135").Value = 8775
$ws.Range("M135").Value = -1778.07693
$ws.Range("N135").Value = -13845

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2473.4482
$ws.Range("I102").Value = 2889.818
$ws.Range("J102").Value = 1164.8572
$ws.Range("K102").Value = 2889.818
$ws.Range("L102").Value = 1164.8572
$ws.Range("M102").Value = -1267.818
$ws.Range("N102").Value = -4408.8572
$ws.Range("H132").Value = 6121.8945
$ws.Range("I132").Value = 5494.5
$ws.Range("J132").Value = 6411.4614
$ws.Range("K132").Value = 16483.5
$ws.Range("L132").Value = 19234.3842
$ws.Range("M132").Value = -13953.5
$ws.Range("N132").Value = -24294.3842
$ws.Range("H138").Value = 56049.5
$ws.Range("J138").Value = 56049.5
$ws.Range("L138").Value = 56049.5
$ws.Range("N138").Value = -66329.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 11116008
$ws.Range("I136").Value = 16668287
$ws.Range("J136").Value = 11448.134
$ws.Range("K136").Value = 50004861
$ws.Range("L136").Value = 34344.402
$ws.Range("M136").Value = -50002311

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4257.7144
$ws.Range("I126").Value = 1634
$ws.Range("J126").Value = 20000
$ws.Range("K126").Value = 4902
$ws.Range("L126").Value = 60000
$ws.Range("M126").Value = -2432
$ws.Range("N126").Value = -64940
$ws.Range("H132").Value = 2115.3125
$ws.Range("I132").Value = 1067.6666
$ws.Range("J132").Value = 3462.2856
$ws.Range("K132").Value = 3202.9998
$ws.Range("L132").Value = 10386.8568
$ws.Range("M132").Value = -672.9998000000001
$ws.Range("N132").Value = -15446.8568
